# add guild data module
# Adds a new "GuilID" property row (row 11) to the "Property" sheet,
# mirroring the existing rows' layout/format (Id, Type, Public, Private,
# Save, View, Index, SaveInterval, RelationValue, Desc).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

$newRow = 11

# Copy formatting (number format / style) from row 2, which already uses
# the plain "text" style (s="1") that the new row should use for its
# string columns (A, B, I, J).
$ws.Range("A2").Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B" + $newRow).PasteSpecial(-4122)

$ws.Range("I2").Copy()
$ws.Range("I" + $newRow).PasteSpecial(-4122)

$ws.Range("J2").Copy()
$ws.Range("J" + $newRow).PasteSpecial(-4122)

# Fill in the values for the new row.
$ws.Range("A" + $newRow).Value = "GuilID"
$ws.Range("B" + $newRow).Value = "object"
$ws.Range("C" + $newRow).Value = $true
$ws.Range("D" + $newRow).Value = $true
$ws.Range("E" + $newRow).Value = $true
$ws.Range("F" + $newRow).Value = $true
$ws.Range("G" + $newRow).Value = 0
$ws.Range("H" + $newRow).Value = 0
$ws.Range("I" + $newRow).Value = "Friend"
$ws.Range("J" + $newRow).Value = "工会ID"

# Update the remembered selection, as recorded in the workbook after the edit.
$ws.Range("E19").Select() | Out-Null
